# Update database: refresh "last updated" headers and quarterly figures
# (monte_carlo / database refresh commit)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Header row 9: publish-date labels shared across columns I/K (9-month) and M (12-month)
$ws.Range("I9").Value = "1402-03-09 (6)"
$ws.Range("K9").Value = "1402-03-09 (6)"
$ws.Range("M9").Value = "1402-03-09 (3)"

# Updated quarterly / cumulative figures
$ws.Range("M12").Value = -11934456
$ws.Range("M13").Value = 5960248
$ws.Range("I14").Value = -78691
$ws.Range("M14").Value = -434524
$ws.Range("I17").Value = 4235096
$ws.Range("M17").Value = 6120641
$ws.Range("M18").Value = -787803
$ws.Range("I20").Value = 3391157
$ws.Range("M20").Value = 5572928
$ws.Range("M21").Value = -1111818
$ws.Range("I22").Value = 2821889
$ws.Range("M22").Value = 4461110
$ws.Range("I24").Value = 2821889
$ws.Range("M24").Value = 4461110
$ws.Range("M25").Value = 312
$ws.Range("M27").Value = 312
